# Generate Report for Archive
#
# The localization status for 025b5ffa-eca5-4a79-b7da-9c2d1fbdd084.md and
# e1474d6f-3107-4f3e-b5ff-caeab908aca6.md moved from "Ready for handoff" to
# "In Translation". Because the rows in each table are kept ordered by file
# name within a status group, this re-sorts rows 7-8 (025b5ffa.md now sorts
# ahead of 83651a89-bb4a-4641-b2b1-c2b339e20237.md, which is also
# "In Translation") while row 9 (e1474d6f.md) and row 10 (fbdf8557.md,
# unaffected) keep their position.

function Set-CellWithHyperlink {
    param(
        $ws,
        [string]$ColLetter,
        [int]$Row,
        [string]$NewValue,
        [bool]$HasHyperlink
    )

    $addr = "$ColLetter$Row"
    $ws.Range($addr).Value2 = $NewValue

    if ($HasHyperlink) {
        $target = '$' + $ColLetter + '$' + $Row
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Address() -eq $target) {
                $h.TextToDisplay = $NewValue
            }
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet: columns A (File Name, hyperlinked), B (zh-cn status),
# C (de-de status), D (Latest Handoff Date)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellWithHyperlink $wsOverview "A" 7  "025b5ffa-eca5-4a79-b7da-9c2d1fbdd084.md" $true
Set-CellWithHyperlink $wsOverview "B" 7  "In Translation"                          $false
Set-CellWithHyperlink $wsOverview "C" 7  "In Translation"                          $false
Set-CellWithHyperlink $wsOverview "D" 7  "2016-18-17 18:18:46"                     $false

Set-CellWithHyperlink $wsOverview "A" 8  "83651a89-bb4a-4641-b2b1-c2b339e20237.md" $true
Set-CellWithHyperlink $wsOverview "B" 8  "In Translation"                          $false
Set-CellWithHyperlink $wsOverview "C" 8  "In Translation"                          $false
Set-CellWithHyperlink $wsOverview "D" 8  "2016-12-17 18:12:49"                     $false

Set-CellWithHyperlink $wsOverview "A" 9  "e1474d6f-3107-4f3e-b5ff-caeab908aca6.md" $true
Set-CellWithHyperlink $wsOverview "B" 9  "In Translation"                          $false
Set-CellWithHyperlink $wsOverview "C" 9  "In Translation"                          $false
Set-CellWithHyperlink $wsOverview "D" 9  "2016-18-17 18:18:46"                     $false

Set-CellWithHyperlink $wsOverview "A" 10 "fbdf8557-9368-407b-a255-6254c559e860.md" $true
Set-CellWithHyperlink $wsOverview "B" 10 "Ready for handoff"                       $false
Set-CellWithHyperlink $wsOverview "C" 10 "Ready for handoff"                       $false
Set-CellWithHyperlink $wsOverview "D" 10 "2016-18-17 18:18:46"                     $false

# ---------------------------------------------------------------------------
# "zh-cn" sheet: column A (Source File Name, hyperlinked), C (Status),
# D (Latest Handoff File, hyperlinked), E (Latest Handoff Datetime)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellWithHyperlink $wsZhCn "A" 7  "025b5ffa-eca5-4a79-b7da-9c2d1fbdd084.md" $true
Set-CellWithHyperlink $wsZhCn "C" 7  "In Translation"                          $false
Set-CellWithHyperlink $wsZhCn "D" 7  "025b5ffa-eca5-4a79-b7da-9c2d1fbdd084.12547876fd5e3fc2394ec3d1bffcf5b3fbf38599.zh-cn.xlf" $true
Set-CellWithHyperlink $wsZhCn "E" 7  "2016-03-17 18:18:42"                     $false

Set-CellWithHyperlink $wsZhCn "A" 8  "83651a89-bb4a-4641-b2b1-c2b339e20237.md" $true
Set-CellWithHyperlink $wsZhCn "C" 8  "In Translation"                          $false
Set-CellWithHyperlink $wsZhCn "D" 8  "83651a89-bb4a-4641-b2b1-c2b339e20237.d4fab33d55e998fea69eb2c73281268644ba286f.zh-cn.xlf" $true
Set-CellWithHyperlink $wsZhCn "E" 8  "2016-03-17 18:12:34"                     $false

Set-CellWithHyperlink $wsZhCn "A" 9  "e1474d6f-3107-4f3e-b5ff-caeab908aca6.md" $true
Set-CellWithHyperlink $wsZhCn "C" 9  "In Translation"                          $false
Set-CellWithHyperlink $wsZhCn "D" 9  "e1474d6f-3107-4f3e-b5ff-caeab908aca6.106c10e42a9a30501a55cbb1515403724d900b77.zh-cn.xlf" $true
Set-CellWithHyperlink $wsZhCn "E" 9  "2016-03-17 18:18:42"                     $false

Set-CellWithHyperlink $wsZhCn "A" 10 "fbdf8557-9368-407b-a255-6254c559e860.md" $true
Set-CellWithHyperlink $wsZhCn "C" 10 "Ready for handoff"                       $false
Set-CellWithHyperlink $wsZhCn "D" 10 "fbdf8557-9368-407b-a255-6254c559e860.59100d35eca6e5d06be1c57423f3ef142a5785ec.zh-cn.xlf" $true
Set-CellWithHyperlink $wsZhCn "E" 10 "2016-03-17 18:18:42"                     $false

# ---------------------------------------------------------------------------
# "de-de" sheet: column A (Source File Name, hyperlinked), C (Status),
# D (Latest Handoff File, hyperlinked), E (Latest Handoff Datetime)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellWithHyperlink $wsDeDe "A" 7  "025b5ffa-eca5-4a79-b7da-9c2d1fbdd084.md" $true
Set-CellWithHyperlink $wsDeDe "C" 7  "In Translation"                          $false
Set-CellWithHyperlink $wsDeDe "D" 7  "025b5ffa-eca5-4a79-b7da-9c2d1fbdd084.12547876fd5e3fc2394ec3d1bffcf5b3fbf38599.de-de.xlf" $true
Set-CellWithHyperlink $wsDeDe "E" 7  "2016-03-17 18:18:46"                     $false

Set-CellWithHyperlink $wsDeDe "A" 8  "83651a89-bb4a-4641-b2b1-c2b339e20237.md" $true
Set-CellWithHyperlink $wsDeDe "C" 8  "In Translation"                          $false
Set-CellWithHyperlink $wsDeDe "D" 8  "83651a89-bb4a-4641-b2b1-c2b339e20237.d4fab33d55e998fea69eb2c73281268644ba286f.de-de.xlf" $true
Set-CellWithHyperlink $wsDeDe "E" 8  "2016-03-17 18:12:49"                     $false

Set-CellWithHyperlink $wsDeDe "A" 9  "e1474d6f-3107-4f3e-b5ff-caeab908aca6.md" $true
Set-CellWithHyperlink $wsDeDe "C" 9  "In Translation"                          $false
Set-CellWithHyperlink $wsDeDe "D" 9  "e1474d6f-3107-4f3e-b5ff-caeab908aca6.106c10e42a9a30501a55cbb1515403724d900b77.de-de.xlf" $true
Set-CellWithHyperlink $wsDeDe "E" 9  "2016-03-17 18:18:46"                     $false

Set-CellWithHyperlink $wsDeDe "A" 10 "fbdf8557-9368-407b-a255-6254c559e860.md" $true
Set-CellWithHyperlink $wsDeDe "C" 10 "Ready for handoff"                       $false
Set-CellWithHyperlink $wsDeDe "D" 10 "fbdf8557-9368-407b-a255-6254c559e860.59100d35eca6e5d06be1c57423f3ef142a5785ec.de-de.xlf" $true
Set-CellWithHyperlink $wsDeDe "E" 10 "2016-03-17 18:18:46"                     $false
